$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LRFR")

# Insert the new "Number of spans" row at row 2, pushing existing rows down.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Fill the new row 2 (description, var, unit) before touching the header row
# so new shared-string entries are created in this order.
$ws.Cells.Item(2, 4).Value = "Number of spans"
$ws.Cells.Item(2, 1).Value = "nSpans"
$ws.Cells.Item(2, 3).Value = "int"
$ws.Cells.Item(2, 2).Value = 2

# New header row: name | value | unit | desc
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "value"
$ws.Range("C1").Value = "unit"
$ws.Range("D1").Value = "desc"

# Remaining data rows (previously rows 2-18, now shifted to rows 3-19):
# var | value | unit | description
$data = @(
    @("L",      1428,     "in",  "Span Length"),
    @("Lb",     324,      "in",  "Max unbraced Length"),
    @("Es",     29000000, "psi", "Steel Modulous"),
    @("Fy",     36000,    "psi", "Steel Yield Strength"),
    @("fc",     4000,     "psi", "Concrete Strength"),
    @("ts",     8.5,      "in",  "Deck Thickness"),
    @("be",     92,       "in",  "Effective Width"),
    @("dh",     4,        "in",  "Haunch Depth"),
    @("dw",     78,       "in",  "Web Depth"),
    @("tw",     0.5625,   "in",  "Web Thickness"),
    @("bf_top", 14,       "in",  "Top Flange Width"),
    @("tf_top", 0.75,     "in",  "Top Flange Thickness"),
    @("bf_bot", 14,       "in",  "Bottom Flange Width"),
    @("tf_bot", 1.5,      "in",  "Bottom Flange Thickness"),
    @("wDL",    109,      "lb",  "NonSuperimposed DL"),
    @("wSDL",   12,       "lb",  "Superimposed DL"),
    @("wSDW",   14,       "lb",  "Wearing DL")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Column widths: A best-fit narrow (var names), D best-fit wide (descriptions)
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(1).ColumnWidth = 6.5
$ws.Columns.Item(4).ColumnWidth = 22.5

# Update selection to match target view state
$ws.Range("J21").Select()
